$wb = $excel.ActiveWorkbook

# Update handoff/handback timestamps on the "zh-cn" report sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 22:52:06"
$wsZh.Range("H2").Value = "2016-03-19 22:52:24"

# Update handoff/handback timestamps on the "de-de" report sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 22:52:09"
$wsDe.Range("H2").Value = "2016-03-19 22:52:30"
